$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 each describe one observation record; the two records trade
# places. Swap every column that actually differs between the rows (Id,
# Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor,
# Ost, Nord). Columns already holding identical data in both rows
# (Valideringsstatus, Rödlistade, Lokalnamn, Noggrannhet, Län, Kommun,
# Provins, Församling, Startdatum, Slutdatum, etc.) are left untouched.
$swapCols = 1, 2, 5, 6, 7, 8, 17, 18   # A, B, E, F, G, H, Q, R

foreach ($col in $swapCols) {
    $cellTop = $ws.Cells.Item(2, $col)
    $cellBottom = $ws.Cells.Item(3, $col)

    $topValue = $cellTop.Value2
    $bottomValue = $cellBottom.Value2

    $cellTop.Value = $bottomValue
    $cellBottom.Value = $topValue
}

# AF (Bestämningsmetod) holds an empty placeholder in only one of the two
# rows; that placeholder moves from row 2 to row 3 along with the record.
$ws.Cells.Item(2, 32).ClearContents()
$ws.Cells.Item(3, 32).Font.Bold = $false
